$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed the three new rows by copying row 5 (value + formatting) so the
# existing style indices (s="1" / s="2") are reused instead of new ones
# being minted in styles.xml.
$ws.Range("A5:F5").Copy($ws.Range("A6:F6"))
$ws.Range("A5:F5").Copy($ws.Range("A7:F7"))
$ws.Range("A5:F5").Copy($ws.Range("A8:F8"))

# Row 6 - com.singleton.strechy / taxi game review
$ws.Range("A6").Value = "com.singleton.strechy"
$ws.Range("B6").Value = "taxi game"
$ws.Range("C6").Value = "redvelvetmichael@gmail.com"
$ws.Range("D6").Value = "veredsnir12@gmail.com"
$ws.Range("E6").Value = "27/5/2019 15:59"
$ws.Range("F6").Value = "Crazy and hard levels but I like it. I can play it all day long. Sweet taxi alabama!!"

# Row 7 - com.hamxa.shaynachim / bitcoin guide review
$ws.Range("A7").Value = "com.hamxa.shaynachim"
$ws.Range("B7").Value = "bitcoin guide"
$ws.Range("C7").Value = "vikicrestina@gmail.com"
$ws.Range("D7").Value = "cristianjohn1222@gmail.com"
$ws.Range("E7").Value = "27/5/2019 15:59"
$ws.Range("F7").Value = "bitcoin guide – great app. Following KISS guidelines – Keep it simple st…"

# Row 8 - com.hamxa.shaynachim / bitcoin guide review
$ws.Range("A8").Value = "com.hamxa.shaynachim"
$ws.Range("B8").Value = "bitcoin guide"
$ws.Range("C8").Value = "bittonnir12@gmail.com"
$ws.Range("D8").Value = "nevilgreen12@gmail.com"
$ws.Range("E8").Value = "27/5/2019 15:59"
$ws.Range("F8").Value = "it is so awesome to get all the bitcoin information in one place. Great development"

# Wire up the mailto hyperlinks for the new email cells.
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:redvelvetmichael@gmail.com", "", "", "redvelvetmichael@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:veredsnir12@gmail.com", "", "", "veredsnir12@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:vikicrestina@gmail.com", "", "", "vikicrestina@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:cristianjohn1222@gmail.com", "", "", "cristianjohn1222@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C8"), "mailto:bittonnir12@gmail.com", "", "", "bittonnir12@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D8"), "mailto:nevilgreen12@gmail.com", "", "", "nevilgreen12@gmail.com")

# Hyperlinks.Add recolors cells with the built-in "Hyperlink" style; restore
# the original plain formatting (style index 2, same as C2/D2) by pasting
# just the formatting back over the top.
$ws.Range("C2:D2").Copy()
$ws.Range("C6:D6").PasteSpecial(-4122)
$ws.Range("C7:D7").PasteSpecial(-4122)
$ws.Range("C8:D8").PasteSpecial(-4122)

# Match the author's final selection (cell F8, the last cell touched).
$ws.Range("F8").Select()
